$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 58

$ws.Cells.Item($newRow, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($newRow, 2).Value = "temperature"
# Leading apostrophe forces this numeric-looking value to be stored as text,
# matching the source row's inline-string "25" (not a number).
$ws.Cells.Item($newRow, 3).Value = "'25"
$ws.Cells.Item($newRow, 4).Value = "N/A"
$ws.Cells.Item($newRow, 5).Value = "N/A"
$ws.Cells.Item($newRow, 6).Value = "N/A"
